$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-23 Wednesday", "2025-07-24 Thursday"),
    @("76÷8=", "37÷2="),
    @("26÷9=", "86÷7="),
    @("64÷9=", "66÷3="),
    @("16÷2=", "65÷4="),
    @("68÷5=", "83÷5="),
    @("74÷5=", "48÷9="),
    @("88÷7=", "66÷4="),
    @("96÷4=", "64÷8="),
    @("69÷8=", "22÷5="),
    @("93÷5=", "54÷2="),
    @("61÷3=", "76÷3="),
    @("86÷4=", "29÷4="),
    @("85÷7=", "20÷2="),
    @("55÷5=", "38÷7="),
    @("45÷3=", "21÷5="),
    @("76÷2=", "90÷9="),
    @("70÷2=", "75÷8="),
    @("21÷7=", "15÷2="),
    @("51÷5=", "12÷8="),
    @("63÷6=", "66÷9="),
    @("34÷6=", "42÷8="),
    @("34÷7=", "76÷2="),
    @("69÷4=", "45÷8="),
    @("81÷9=", "29÷3="),
    @("63÷8=", "81÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
